# Append 54 new daily-sentiment rows (rows 3525-3578) after the existing
# data which currently ends at row 3524 (A1:B3524).
# Column A = serial date value (formatted like the existing date column),
# Column B = sentiment value (0 for all new rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 3525

$dates = @(
    45838, 45839, 45840, 45841, 45845,
    45846, 45847, 45848, 45849, 45852,
    45853, 45854, 45855, 45856, 45859,
    45860, 45861, 45862, 45863, 45866,
    45867, 45868, 45869, 45870, 45873,
    45874, 45875, 45876, 45877, 45880,
    45881, 45882, 45883, 45884, 45887,
    45888, 45889, 45890, 45891, 45894,
    45895, 45896, 45897, 45898, 45902,
    45903, 45904, 45905, 45908, 45909,
    45910, 45911, 45912, 45915
)

$sourceStyle = $ws.Range("A" + ($startRow - 1))
$targetFormat = $sourceStyle.NumberFormat

for ($i = 0; $i -lt $dates.Count; $i++) {
    $r = $startRow + $i

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $dates[$i]
    $cellA.NumberFormat = $targetFormat

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = 0
}

Write-Host "Added rows" $startRow "to" ($startRow + $dates.Count - 1)
